$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Tuesday office-hours line (Word paragraph #3 before any edits):
#   "Tuesday      2:30 to 4:00"  ->  "Tuesday      10:00 - 11:00"
# and a brand-new line is added right after it:
#   "Tuesday      12:15 - 1:00"
# ---------------------------------------------------------------------

$tuesRange = $d.Paragraphs(3).Range.Duplicate
$tuesRange.Find.ClearFormatting()
$foundTues = $tuesRange.Find.Execute(
    "Tuesday      2:30 to 4:00", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tuesday      10:00 - 11:00", 2)
if (-not $foundTues) {
    throw "Could not find the Tuesday office-hours line to update."
}

$tuesEnd = $d.Paragraphs(3).Range.Duplicate
$tuesEnd.Collapse(0)
$tuesEnd.InsertParagraphAfter()

$tuesNew = $d.Paragraphs(4)
$tuesNew.Range.Font.Size = 30
$tuesNew.Range.Text = "Tuesday      12:15 - 1:00"

# ---------------------------------------------------------------------
# Thursday office-hours line (Word paragraph #5 once the Tuesday insert
# above has shifted everything below it down by one paragraph):
#   "Thursday     2:30 to 4:00"  ->  "Thursday     10:00 -11:00"
# and a brand-new line is added right after it:
#   "Thursday     12:15 -1:00"
# ---------------------------------------------------------------------

$thuRange = $d.Paragraphs(5).Range.Duplicate
$thuRange.Find.ClearFormatting()
$foundThu = $thuRange.Find.Execute(
    "Thursday     2:30 to 4:00", $true, $false, $false, $false, $false,
    $true, 1, $false, "Thursday     10:00 -11:00", 2)
if (-not $foundThu) {
    throw "Could not find the Thursday office-hours line to update."
}

$thuEnd = $d.Paragraphs(5).Range.Duplicate
$thuEnd.Collapse(0)
$thuEnd.InsertParagraphAfter()

$thuNew = $d.Paragraphs(6)
$thuNew.Range.Font.Size = 30
$thuNew.Range.Text = "Thursday     12:15 -1:00"
